$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.361.21'
$ws.Range('E2').Value = '  -2.12%  '

$ws.Range('D3').Value = '1.708.35'
$ws.Range('E3').Value = '  -1.94%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.003'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.07%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '223.88'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.93%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5329'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.80%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.003'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.03%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2660'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -3.83%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06590'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -2.10%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '20.84'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -3.80%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07640'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.80%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.568'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -2.94%  '

$ws.Range('D13').Value = '1.704.93'
$ws.Range('E13').Value = '  -1.98%  '

$ws.Range('D14').Value = '1.943.17'
$ws.Range('E14').Value = '  -2.03%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.5727'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -4.06%  '

$ws.Range('D16').Value = '0.0₅8173'
$ws.Range('E16').Value = '  -2.54%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '67.71'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.68%  '

$ws.Range('D18').Value = '27.325.64'
$ws.Range('E18').Value = '  -2.31%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '216.27'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -3.52%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.003'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.00%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.669'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -3.59%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.44'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -4.25%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.975'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -4.23%  '

$ws.Range('E24').Value = '  -0.09%  '

$ws.Range('E25').Value = '  +5.79%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '141.84'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -3.12%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.1215'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.38%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.268'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -2.49%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '16.32'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -5.41%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.05424'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -3.76%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.293'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.73%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.504'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -5.46%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.427'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.64%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.645'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -2.08%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.875'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.66%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.9490'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -3.26%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.407'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.53%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.5864'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.77%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01631'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.13%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.858'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.38%  '

$ws.Range('D41').Value = '1.044.34'
$ws.Range('E41').Value = '  -0.36%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.8432'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.79%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.003'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.03%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '100.85'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.25%  '

$ws.Range('D45').Value = '1.850.34'
$ws.Range('E45').Value = '  -2.02%  '

$ws.Range('E46').Value = '  -2.32%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '58.02'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -3.22%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.4501'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.43%  '

$ws.Range('E49').Value = '  -0.20%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '8.093'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.45%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.05243'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.50%  '
